$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 1697.5  # ALC!H18 was 1500
$ws.Cells.Item(18, 9).Value = 1597  # ALC!I18 was 1500
$ws.Cells.Item(18, 10).Value = 2200  # ALC!J18 was 0
$ws.Cells.Item(18, 11).Value = 1597  # ALC!K18 was 1500
$ws.Cells.Item(18, 12).Value = 2200  # ALC!L18 was 0
$ws.Cells.Item(18, 13).Value = -1313  # ALC!M18 was -1216
$ws.Cells.Item(18, 14).Value = -2768  # ALC!N18 was None
$ws.Cells.Item(40, 8).Value = 45456544  # ALC!H40 was 41668670
$ws.Cells.Item(40, 9).Value = 1250.25  # ALC!I40 was 1433.3334
$ws.Cells.Item(40, 10).Value = 55557724  # ALC!J40 was 47621132
$ws.Cells.Item(40, 11).Value = 1250.25  # ALC!K40 was 1433.3334
$ws.Cells.Item(40, 12).Value = 55557724  # ALC!L40 was 47621132
$ws.Cells.Item(40, 13).Value = -1075.25  # ALC!M40 was -1258.3334
$ws.Cells.Item(40, 14).Value = -55558074  # ALC!N40 was -47621482
$ws.Cells.Item(88, 8).Value = 2647793.2  # ALC!H88 was 3969003.8
$ws.Cells.Item(88, 9).Value = 10000  # ALC!I88 was 1000
$ws.Cells.Item(88, 10).Value = 3175351.8  # ALC!J88 was 5291671.5
$ws.Cells.Item(88, 11).Value = 10000  # ALC!K88 was 1000
$ws.Cells.Item(88, 12).Value = 3175351.8  # ALC!L88 was 5291671.5
$ws.Cells.Item(88, 13).Value = -9594  # ALC!M88 was -594
$ws.Cells.Item(88, 14).Value = -3176163.8  # ALC!N88 was -5292483.5
$ws.Cells.Item(91, 8).Value = 2647793.2  # ALC!H91 was 3969003.8
$ws.Cells.Item(91, 9).Value = 10000  # ALC!I91 was 1000
$ws.Cells.Item(91, 10).Value = 3175351.8  # ALC!J91 was 5291671.5
$ws.Cells.Item(91, 11).Value = 10000  # ALC!K91 was 1000
$ws.Cells.Item(91, 12).Value = 3175351.8  # ALC!L91 was 5291671.5
$ws.Cells.Item(91, 13).Value = -8596  # ALC!M91 was 404
$ws.Cells.Item(91, 14).Value = -3178159.8  # ALC!N91 was -5294479.5
$ws.Cells.Item(106, 8).Value = 3157  # ALC!H106 was 2296.5
$ws.Cells.Item(106, 9).Value = 3274.75  # ALC!I106 was 2134.1538
$ws.Cells.Item(106, 11).Value = 3274.75  # ALC!K106 was 2134.1538
$ws.Cells.Item(106, 13).Value = -2643.75  # ALC!M106 was -1503.1538
$ws.Cells.Item(134, 8).Value = 55823.53  # ALC!H134 was 67161.336
$ws.Cells.Item(134, 10).Value = 55823.53  # ALC!J134 was 67161.336
$ws.Cells.Item(134, 12).Value = 55823.53  # ALC!L134 was 67161.336
$ws.Cells.Item(134, 14).Value = -65963.53  # ALC!N134 was -77301.336
$ws.Cells.Item(135, 8).Value = 10001  # ALC!H135 was 0
$ws.Cells.Item(135, 9).Value = 10001  # ALC!I135 was 0
$ws.Cells.Item(135, 11).Value = 90009  # ALC!K135 was 0
$ws.Cells.Item(135, 13).Value = -87474  # ALC!M135 was None
$ws.Cells.Item(141, 8).Value = 1721.0714  # ALC!H141 was 1421.0714
$ws.Cells.Item(141, 9).Value = 1439  # ALC!I141 was 1476.1538
$ws.Cells.Item(141, 10).Value = 1877.7778  # ALC!J141 was 705
$ws.Cells.Item(141, 11).Value = 4317  # ALC!K141 was 4428.4614
$ws.Cells.Item(141, 12).Value = 5633.3334  # ALC!L141 was 2115
$ws.Cells.Item(141, 13).Value = 863  # ALC!M141 was 751.5385999999999
$ws.Cells.Item(141, 14).Value = -15993.3334  # ALC!N141 was -12475

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 17797.602  # ARM!H32 was 17723.215
$ws.Cells.Item(32, 9).Value = 15168.855  # ARM!I32 was 15770.453
$ws.Cells.Item(32, 10).Value = 96660  # ARM!J32 was 41714.285
$ws.Cells.Item(32, 11).Value = 15168.855  # ARM!K32 was 15770.453
$ws.Cells.Item(32, 12).Value = 96660  # ARM!L32 was 41714.285
$ws.Cells.Item(32, 13).Value = -14881.855  # ARM!M32 was -15483.453
$ws.Cells.Item(32, 14).Value = -97234  # ARM!N32 was -42288.285
$ws.Cells.Item(45, 8).Value = 8696.883  # ARM!H45 was 3451.9167
$ws.Cells.Item(45, 9).Value = 10078.917  # ARM!I45 was 2170.5
$ws.Cells.Item(45, 10).Value = 5380  # ARM!J45 was 4733.3335
$ws.Cells.Item(45, 11).Value = 10078.917  # ARM!K45 was 2170.5
$ws.Cells.Item(45, 12).Value = 5380  # ARM!L45 was 4733.3335
$ws.Cells.Item(45, 13).Value = -9701.916999999999  # ARM!M45 was -1793.5
$ws.Cells.Item(45, 14).Value = -6134  # ARM!N45 was -5487.3335
$ws.Cells.Item(74, 8).Value = 808.4902  # ARM!H74 was 937.08887
$ws.Cells.Item(74, 9).Value = 992.53845  # ARM!I74 was 1194.0952
$ws.Cells.Item(74, 10).Value = 617.08  # ARM!J74 was 712.2083
$ws.Cells.Item(74, 11).Value = 992.53845  # ARM!K74 was 1194.0952
$ws.Cells.Item(74, 12).Value = 617.08  # ARM!L74 was 712.2083
$ws.Cells.Item(74, 13).Value = -118.53845  # ARM!M74 was -320.0952
$ws.Cells.Item(74, 14).Value = -2365.08  # ARM!N74 was -2460.2083
$ws.Cells.Item(77, 8).Value = 808.4902  # ARM!H77 was 937.08887
$ws.Cells.Item(77, 9).Value = 992.53845  # ARM!I77 was 1194.0952
$ws.Cells.Item(77, 10).Value = 617.08  # ARM!J77 was 712.2083
$ws.Cells.Item(77, 11).Value = 4962.69225  # ARM!K77 was 5970.476
$ws.Cells.Item(77, 12).Value = 3085.4  # ARM!L77 was 3561.0415
$ws.Cells.Item(77, 13).Value = -594.6922500000001  # ARM!M77 was -1602.476
$ws.Cells.Item(77, 14).Value = -11821.4  # ARM!N77 was -12297.0415
$ws.Cells.Item(110, 8).Value = 738.5789  # ARM!H110 was 738.8
$ws.Cells.Item(110, 9).Value = 676.26666  # ARM!I110 was 696.5294
$ws.Cells.Item(110, 10).Value = 972.25  # ARM!J110 was 978.3333
$ws.Cells.Item(110, 11).Value = 676.26666  # ARM!K110 was 696.5294
$ws.Cells.Item(110, 12).Value = 972.25  # ARM!L110 was 978.3333
$ws.Cells.Item(110, 13).Value = 1368.73334  # ARM!M110 was 1348.4706
$ws.Cells.Item(110, 14).Value = -5062.25  # ARM!N110 was -5068.3333
$ws.Cells.Item(139, 8).Value = 49857.5  # ARM!H139 was 55715
$ws.Cells.Item(139, 10).Value = 49857.5  # ARM!J139 was 55715
$ws.Cells.Item(139, 12).Value = 49857.5  # ARM!L139 was 55715
$ws.Cells.Item(139, 14).Value = -60137.5  # ARM!N139 was -65995

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 4083.6365  # BSM!H105 was 15666.667
$ws.Cells.Item(105, 9).Value = 4192  # BSM!I105 was 8500
$ws.Cells.Item(105, 10).Value = 3000  # BSM!J105 was 30000
$ws.Cells.Item(105, 11).Value = 4192  # BSM!K105 was 8500
$ws.Cells.Item(105, 12).Value = 3000  # BSM!L105 was 30000
$ws.Cells.Item(105, 13).Value = -2445  # BSM!M105 was -6753
$ws.Cells.Item(105, 14).Value = -6494  # BSM!N105 was -33494

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 163.66667  # CRP!H7 was 169.76471
$ws.Cells.Item(7, 9).Value = 146.71428  # CRP!I7 was 153.38461
$ws.Cells.Item(7, 11).Value = 146.71428  # CRP!K7 was 153.38461
$ws.Cells.Item(7, 13).Value = -33.71428  # CRP!M7 was -40.38461000000001
$ws.Cells.Item(16, 8).Value = 0  # CRP!H16 was 894
$ws.Cells.Item(16, 9).Value = 0  # CRP!I16 was 564.5
$ws.Cells.Item(16, 10).Value = 0  # CRP!J16 was 1333.3334
$ws.Cells.Item(16, 11).Value = 0  # CRP!K16 was 564.5
$ws.Cells.Item(16, 12).Value = 0  # CRP!L16 was 1333.3334
$ws.Cells.Item(16, 13).ClearContents()  # CRP!M16 was -277.5
$ws.Cells.Item(16, 14).ClearContents()  # CRP!N16 was -1907.3334
$ws.Cells.Item(31, 8).Value = 2063.4333  # CRP!H31 was 2121.5173
$ws.Cells.Item(31, 9).Value = 909.8889  # CRP!I31 was 1005
$ws.Cells.Item(31, 10).Value = 2557.8096  # CRP!J31 was 2354.125
$ws.Cells.Item(31, 11).Value = 909.8889  # CRP!K31 was 1005
$ws.Cells.Item(31, 12).Value = 2557.8096  # CRP!L31 was 2354.125
$ws.Cells.Item(31, 13).Value = -614.8889  # CRP!M31 was -710
$ws.Cells.Item(31, 14).Value = -3147.8096  # CRP!N31 was -2944.125
$ws.Cells.Item(34, 8).Value = 2063.4333  # CRP!H34 was 2121.5173
$ws.Cells.Item(34, 9).Value = 909.8889  # CRP!I34 was 1005
$ws.Cells.Item(34, 10).Value = 2557.8096  # CRP!J34 was 2354.125
$ws.Cells.Item(34, 11).Value = 909.8889  # CRP!K34 was 1005
$ws.Cells.Item(34, 12).Value = 2557.8096  # CRP!L34 was 2354.125
$ws.Cells.Item(34, 13).Value = -707.8889  # CRP!M34 was -803
$ws.Cells.Item(34, 14).Value = -2961.8096  # CRP!N34 was -2758.125
$ws.Cells.Item(60, 8).Value = 4850  # CRP!H60 was 10000
$ws.Cells.Item(60, 9).Value = 4850  # CRP!I60 was 10000
$ws.Cells.Item(60, 11).Value = 4850  # CRP!K60 was 10000
$ws.Cells.Item(60, 13).Value = -4339  # CRP!M60 was -9489
$ws.Cells.Item(68, 8).Value = 21225  # CRP!H68 was 22771.072
$ws.Cells.Item(68, 10).Value = 21225  # CRP!J68 was 22771.072
$ws.Cells.Item(68, 12).Value = 21225  # CRP!L68 was 22771.072
$ws.Cells.Item(68, 14).Value = -22723  # CRP!N68 was -24269.072
$ws.Cells.Item(71, 8).Value = 21225  # CRP!H71 was 22771.072
$ws.Cells.Item(71, 10).Value = 21225  # CRP!J71 was 22771.072
$ws.Cells.Item(71, 12).Value = 63675  # CRP!L71 was 68313.216
$ws.Cells.Item(71, 14).Value = -71163  # CRP!N71 was -75801.216
$ws.Cells.Item(74, 8).Value = 0  # CRP!H74 was 24907
$ws.Cells.Item(74, 10).Value = 0  # CRP!J74 was 24907
$ws.Cells.Item(74, 12).Value = 0  # CRP!L74 was 24907
$ws.Cells.Item(74, 14).ClearContents()  # CRP!N74 was -26655
$ws.Cells.Item(77, 8).Value = 0  # CRP!H77 was 24907
$ws.Cells.Item(77, 10).Value = 0  # CRP!J77 was 24907
$ws.Cells.Item(77, 12).Value = 0  # CRP!L77 was 74721
$ws.Cells.Item(77, 14).ClearContents()  # CRP!N77 was -83457
$ws.Cells.Item(105, 8).Value = 1795.2858  # CRP!H105 was 2078.2273
$ws.Cells.Item(105, 9).Value = 1780  # CRP!I105 was 1986.4286
$ws.Cells.Item(105, 10).Value = 1820.125  # CRP!J105 was 2238.875
$ws.Cells.Item(105, 11).Value = 1780  # CRP!K105 was 1986.4286
$ws.Cells.Item(105, 12).Value = 1820.125  # CRP!L105 was 2238.875
$ws.Cells.Item(105, 13).Value = -33  # CRP!M105 was -239.4286
$ws.Cells.Item(105, 14).Value = -5314.125  # CRP!N105 was -5732.875
$ws.Cells.Item(113, 8).Value = 0  # CRP!H113 was 894
$ws.Cells.Item(113, 9).Value = 0  # CRP!I113 was 564.5
$ws.Cells.Item(113, 10).Value = 0  # CRP!J113 was 1333.3334
$ws.Cells.Item(113, 11).Value = 0  # CRP!K113 was 564.5
$ws.Cells.Item(113, 12).Value = 0  # CRP!L113 was 1333.3334
$ws.Cells.Item(113, 13).ClearContents()  # CRP!M113 was 1605.5
$ws.Cells.Item(113, 14).ClearContents()  # CRP!N113 was -5673.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 1251.0597  # CUL!H68 was 1547.1464
$ws.Cells.Item(68, 9).Value = 762  # CUL!I68 was 0
$ws.Cells.Item(68, 10).Value = 1524.0233  # CUL!J68 was 1547.1464
$ws.Cells.Item(68, 11).Value = 2286  # CUL!K68 was 0
$ws.Cells.Item(68, 12).Value = 4572.0699  # CUL!L68 was 4641.439200000001
$ws.Cells.Item(68, 13).Value = -1475  # CUL!M68 was None
$ws.Cells.Item(68, 14).Value = -6194.0699  # CUL!N68 was -6263.439200000001
$ws.Cells.Item(71, 8).Value = 1251.0597  # CUL!H71 was 1547.1464
$ws.Cells.Item(71, 9).Value = 762  # CUL!I71 was 0
$ws.Cells.Item(71, 10).Value = 1524.0233  # CUL!J71 was 1547.1464
$ws.Cells.Item(71, 11).Value = 6858  # CUL!K71 was 0
$ws.Cells.Item(71, 12).Value = 13716.2097  # CUL!L71 was 13924.3176
$ws.Cells.Item(71, 13).Value = -2802  # CUL!M71 was None
$ws.Cells.Item(71, 14).Value = -21828.2097  # CUL!N71 was -22036.3176
$ws.Cells.Item(113, 8).Value = 167219.17  # CUL!H113 was 175997.38
$ws.Cells.Item(113, 10).Value = 175999.12  # CUL!J113 was 179133.03
$ws.Cells.Item(113, 12).Value = 527997.36  # CUL!L113 was 537399.09
$ws.Cells.Item(113, 14).Value = -532337.36  # CUL!N113 was -541739.09
$ws.Cells.Item(117, 8).Value = 3130.5  # CUL!H117 was 4574.3335
$ws.Cells.Item(117, 9).Value = 995.6667  # CUL!I117 was 0
$ws.Cells.Item(117, 10).Value = 4045.4285  # CUL!J117 was 4574.3335
$ws.Cells.Item(117, 11).Value = 2987.0001  # CUL!K117 was 0
$ws.Cells.Item(117, 12).Value = 12136.2855  # CUL!L117 was 13723.0005
$ws.Cells.Item(117, 13).Value = 454.9998999999998  # CUL!M117 was None
$ws.Cells.Item(117, 14).Value = -19020.2855  # CUL!N117 was -20607.0005
$ws.Cells.Item(132, 8).Value = 659.4  # CUL!H132 was 792.8570999999999
$ws.Cells.Item(132, 9).Value = 489.81818  # CUL!I132 was 650
$ws.Cells.Item(132, 10).Value = 866.6667  # CUL!J132 was 900
$ws.Cells.Item(132, 11).Value = 4408.36362  # CUL!K132 was 5850
$ws.Cells.Item(132, 12).Value = 7800.0003  # CUL!L132 was 8100
$ws.Cells.Item(132, 13).Value = -1878.36362  # CUL!M132 was -3320
$ws.Cells.Item(132, 14).Value = -12860.0003  # CUL!N132 was -13160
$ws.Cells.Item(133, 8).Value = 600  # CUL!H133 was 0
$ws.Cells.Item(133, 9).Value = 600  # CUL!I133 was 0
$ws.Cells.Item(133, 11).Value = 1800  # CUL!K133 was 0
$ws.Cells.Item(133, 13).Value = 3260  # CUL!M133 was None

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1683.5  # GSM!H102 was 932.5
$ws.Cells.Item(102, 9).Value = 812.75  # GSM!I102 was 735
$ws.Cells.Item(102, 10).Value = 3425  # GSM!J102 was 3500
$ws.Cells.Item(102, 11).Value = 812.75  # GSM!K102 was 735
$ws.Cells.Item(102, 12).Value = 3425  # GSM!L102 was 3500
$ws.Cells.Item(102, 13).Value = 809.25  # GSM!M102 was 887
$ws.Cells.Item(102, 14).Value = -6669  # GSM!N102 was -6744
$ws.Cells.Item(135, 8).Value = 40000  # GSM!H135 was 29785
$ws.Cells.Item(135, 10).Value = 40000  # GSM!J135 was 29785
$ws.Cells.Item(135, 12).Value = 40000  # GSM!L135 was 29785
$ws.Cells.Item(135, 14).Value = -50140  # GSM!N135 was -39925

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 290580.3  # LTW!H40 was 237990.2
$ws.Cells.Item(40, 9).Value = 298768.3  # LTW!I40 was 275262.75
$ws.Cells.Item(40, 10).Value = 4000  # LTW!J40 was 1930.8334
$ws.Cells.Item(40, 11).Value = 298768.3  # LTW!K40 was 275262.75
$ws.Cells.Item(40, 12).Value = 4000  # LTW!L40 was 1930.8334
$ws.Cells.Item(40, 13).Value = -298632.3  # LTW!M40 was -275126.75
$ws.Cells.Item(40, 14).Value = -4272  # LTW!N40 was -2202.8334
$ws.Cells.Item(46, 8).Value = 453.07693  # LTW!H46 was 381.90475
$ws.Cells.Item(46, 9).Value = 415  # LTW!I46 was 324.16666
$ws.Cells.Item(46, 10).Value = 485.7143  # LTW!J46 was 458.8889
$ws.Cells.Item(46, 11).Value = 415  # LTW!K46 was 324.16666
$ws.Cells.Item(46, 12).Value = 485.7143  # LTW!L46 was 458.8889
$ws.Cells.Item(46, 13).Value = -227  # LTW!M46 was -136.16666
$ws.Cells.Item(46, 14).Value = -861.7143  # LTW!N46 was -834.8888999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1437.2037  # WVR!H132 was 1452.2587
$ws.Cells.Item(132, 9).Value = 1031.4667  # WVR!I132 was 1137.0667
$ws.Cells.Item(132, 10).Value = 3465.889  # WVR!J132 was 2543.3076
$ws.Cells.Item(132, 11).Value = 3094.4001  # WVR!K132 was 3411.2001
$ws.Cells.Item(132, 12).Value = 10397.667  # WVR!L132 was 7629.9228
$ws.Cells.Item(132, 13).Value = -564.4000999999998  # WVR!M132 was -881.2001
$ws.Cells.Item(132, 14).Value = -15457.667  # WVR!N132 was -12689.9228
